# Dharamshala admission letter: remove the second (co-signer) signature
# block so that only Nitesh Sharma's signature / title / phone number
# remain, per "updated D/shala campus info".

$d = $word.ActiveDocument

# 1. Remove the second signature image (image6.jpg) that sat next to
#    Nitesh Sharma's signature image.
$d.InlineShapes.Item(2).Delete()

# 2. Drop "Rahit Roy" from the printed-name line, keeping Nitesh Sharma's
#    name and the separating whitespace.
$d.Content.Find.Execute( `
    "Nitesh Sharma                                                                                                    Rahit Roy", `
    $true, $false, $false, $false, $false, $true, 1, $false, `
    "Nitesh Sharma                                                                                                    ", 2)

# 3. Drop the "Boys Campus Facility In-Charge" job title, keeping the
#    leading spaces that followed the tab stop.
$d.Content.Find.Execute( `
    "      Boys Campus Facility In-Charge", `
    $true, $false, $false, $false, $false, $true, 1, $false, `
    "      ", 2)

# 4. Drop the trailing phone number, keeping the leading spaces.
$d.Content.Find.Execute( `
    "                                                +91-9354978726", `
    $true, $false, $false, $false, $false, $true, 1, $false, `
    "                                                ", 2)
